$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with full content swap (coin identity moved between adjacent rows) ---
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.061'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.65%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9022'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.19%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.802'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.72%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.23'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.58%  '

# --- Price (D) column updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.240.32'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.855.91'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.00'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4788'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2808'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06489'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.865.33'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07351'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.27'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.117'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.23'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6485'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.201.41'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.22'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.0000'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007629'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '224.97'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.094.44'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.287'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.073'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.250'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.24'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.45'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.921'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.443'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09203'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.242'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.962'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05019'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7389'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.152'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.689'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01832'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.612'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.953'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.36'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4259'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1321'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.376'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.543'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.39'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05664'

# --- Volume(1h) (E) column updates ---
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.25%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.51%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.36%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.32%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.75%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.35%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.69%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.88%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.84%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +17.08%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.01%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.87%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.58%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.10%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.59%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.39%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.16%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.52%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.23%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.08%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.35%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.09%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.77%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.17%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +9.50%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.88%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.87%  '
